$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-code (column C) and codeforiati:group-name (column D)
# columns - including their header cells - have swapped places. Swap the
# whole column contents (header + all data rows) in place.
$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $codeCell = $ws.Cells.Item($r, 3)
    $nameCell = $ws.Cells.Item($r, 4)

    $codeVal = $codeCell.Value2
    $nameVal = $nameCell.Value2

    $codeCell.Value = $nameVal
    $nameCell.Value = $codeVal
}
